# "Remove form_id from basic forms"
#
# The "settings" sheet has columns: form_title | form_id | version | style | namespaces
# This removes the form_id column (column B) entirely, which shifts
# version/style/namespaces left by one column (version -> B, style -> C, namespaces -> D),
# along with their row-2 values and header-row cell comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# The column headers live on row 1 and carry cell comments describing each
# setting. Column deletion via COM does not automatically re-home comments
# to the shifted cells, so capture their text up front and rebuild them
# after the column is removed.
$commentTexts = @{}
for ($col = 1; $col -le 5; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $commentTexts[$col] = $cell.Comment.Text()
    }
}

# Remove all of the existing comments in row 1 before shifting columns
# around, so none are left "stranded" on the wrong cell afterwards.
for ($col = 1; $col -le 5; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete()
    }
}

# Delete column B (form_id). This shifts columns C, D, E (version, style,
# namespaces) left into B, C, D, taking their row 1/2 contents with them.
$ws.Columns.Item(2).Delete()

# Re-create the comments on their new (shifted) cells, skipping the one
# that belonged to the now-deleted form_id column.
$ws.Cells.Item(1, 1).AddComment($commentTexts[1])
$ws.Cells.Item(1, 2).AddComment($commentTexts[3])
$ws.Cells.Item(1, 3).AddComment($commentTexts[4])
$ws.Cells.Item(1, 4).AddComment($commentTexts[5])
